$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).ColumnWidth = 15.65
$ws.Columns.Item(2).ColumnWidth = 14.65
$ws.Range("A1").Value = -0.22510864227085392
$ws.Range("B1").Value = 0.22495396144379498
$ws.Range("A2").Value = -0.17148038099303342
$ws.Range("B2").Value = 0.17101768648356241
$ws.Range("A3").Value = -0.089263900635364379
$ws.Range("B3").Value = 0.08901114280356559
$ws.Range("A4").Value = -0.081011142869305885
$ws.Range("B4").Value = 0.08059621719089094
$ws.Range("A5").Value = -0.077596217228030895
$ws.Range("B5").Value = 0.07618953800106798
$ws.Range("A6").Value = 0.025119133493104684
$ws.Range("B6").Value = -0.02517037660928878
$ws.Range("A7").Value = 0.035170376519428448
$ws.Range("B7").Value = -0.035189275364993833
$ws.Range("A8").Value = -0.00012255912986969264
$ws.Range("B8").Value = 0.00009324093813090073
$ws.Range("A9").Value = 0.0019067590257253464
$ws.Range("B9").Value = -0.0019236989372219959
$ws.Range("A10").Value = 0.0039236989012270129
$ws.Range("B10").Value = -0.0039234520045994259
$ws.Range("A11").Value = 0.0069234519622307644
$ws.Range("B11").Value = -0.0069235687576174598
$ws.Range("A12").Value = 0.0104235687122336
$ws.Range("B12").Value = -0.010432093996584069
$ws.Range("A13").Value = 0.013932093952814917
$ws.Range("B13").Value = -0.013941258453879612
$ws.Range("A14").Value = 0.021941258382046414
$ws.Range("B14").Value = -0.021965459540300714
$ws.Range("A15").Value = 0.022965459514595388
$ws.Range("B15").Value = -0.023000239133644129
$ws.Range("A16").Value = -0.0060324273063834255
$ws.Range("B16").Value = 0.0060031787155572225
$ws.Range("A17").Value = -0.0040031787467293967
$ws.Range("B17").Value = 0.0039999999557558397
$ws.Range("A18").Value = -0.016103756730686314
$ws.Range("B18").Value = 0.016091150696329493
$ws.Range("A19").Value = -0.012091150723576582
$ws.Range("B19").Value = 0.012016438663065898
$ws.Range("A20").Value = -0.008016438692338923
$ws.Range("B20").Value = 0.0080055903315585653
$ws.Range("A21").Value = -0.0040055903611388999
$ws.Range("B21").Value = 0.0039999999701594291
$ws.Range("A22").Value = -0.077753785884935311
$ws.Range("B22").Value = 0.077458274744609312
$ws.Range("A23").Value = -0.072458274791543431
$ws.Range("B23").Value = 0.071884115817107741
$ws.Range("A24").Value = -0.051884115968025668
$ws.Range("B24").Value = 0.051675097983191698
$ws.Range("A25").Value = -0.097194325611425825
$ws.Range("B25").Value = 0.09707321465194596
$ws.Range("A26").Value = -0.094573214696371366
$ws.Range("B26").Value = 0.094416867038502161
$ws.Range("A27").Value = -0.091916867085744425
$ws.Range("B27").Value = 0.090988622722201651
$ws.Range("A28").Value = -0.088988622779091031
$ws.Range("B28").Value = 0.088353173622129688
$ws.Range("A29").Value = -0.081353173720598804
$ws.Range("B29").Value = 0.08116910508842512
$ws.Range("A30").Value = -0.021169105531039456
$ws.Range("B30").Value = 0.021022649404528249
$ws.Range("A31").Value = -0.01402264951078358
$ws.Range("B31").Value = 0.014001007428086965
$ws.Range("A32").Value = -0.0040010075540557466
$ws.Range("B32").Value = 0.0039999999123025987
